$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4918.4
$ws.Range("I51").Value = 5059.5713
$ws.Range("K51").Value = 5059.5713
$ws.Range("M51").Value = -4575.5713
$ws.Range("H62").Value = 3942.5715
$ws.Range("I62").Value = 3333.3333
$ws.Range("K62").Value = 3333.3333
$ws.Range("M62").Value = -2709.3333
$ws.Range("H65").Value = 3942.5715
$ws.Range("I65").Value = 3333.3333
$ws.Range("K65").Value = 16666.6665
$ws.Range("M65").Value = -13546.6665
$ws.Range("H86").Value = 3235000
$ws.Range("I86").Value = 4619009
$ws.Range("K86").Value = 4619009
$ws.Range("M86").Value = -4617886
$ws.Range("H89").Value = 3235000
$ws.Range("I89").Value = 4619009
$ws.Range("K89").Value = 23095045
$ws.Range("M89").Value = -23089429
$ws.Range("H100").Value = 1009.3333
$ws.Range("J100").Value = 1499
$ws.Range("L100").Value = 1499
$ws.Range("N100").Value = -2581
$ws.Range("H125").Value = 3395.889
$ws.Range("I125").Value = 2032
$ws.Range("J125").Value = 3566.375
$ws.Range("K125").Value = 18288
$ws.Range("L125").Value = 32097.375
$ws.Range("M125").Value = -15828
$ws.Range("N125").Value = -37017.375
$ws.Range("H137").Value = 8776406
$ws.Range("I137").Value = 1730.8667
$ws.Range("J137").Value = 14499020
$ws.Range("K137").Value = 5192.6001
$ws.Range("L137").Value = 43497060
$ws.Range("M137").Value = -2642.6001
$ws.Range("N137").Value = -43502160
$ws.Range("H138").Value = 3199.818
$ws.Range("I138").Value = 1766.25
$ws.Range("J138").Value = 3518.389
$ws.Range("K138").Value = 5298.75
$ws.Range("L138").Value = 10555.167
$ws.Range("M138").Value = -158.75
$ws.Range("N138").Value = -20835.167
$ws.Range("H141").Value = 12706.8
$ws.Range("I141").Value = 12706.8
$ws.Range("K141").Value = 38120.39999999999
$ws.Range("M141").Value = -32940.39999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4965.6665
$ws.Range("I45").Value = 4990
$ws.Range("J45").Value = 4960.8
$ws.Range("K45").Value = 4990
$ws.Range("L45").Value = 4960.8
$ws.Range("M45").Value = -4613
$ws.Range("N45").Value = -5714.8
$ws.Range("H61").Value = 10731.4375
$ws.Range("I61").Value = 11313.6
$ws.Range("K61").Value = 11313.6
$ws.Range("M61").Value = -11101.6
$ws.Range("H74").Value = 11365553
$ws.Range("I74").Value = 16668138
$ws.Range("J74").Value = 2870.5715
$ws.Range("K74").Value = 16668138
$ws.Range("L74").Value = 2870.5715
$ws.Range("M74").Value = -16667264
$ws.Range("N74").Value = -4618.5715
$ws.Range("H77").Value = 11365553
$ws.Range("I77").Value = 16668138
$ws.Range("J77").Value = 2870.5715
$ws.Range("K77").Value = 83340690
$ws.Range("L77").Value = 14352.8575
$ws.Range("M77").Value = -83336322
$ws.Range("N77").Value = -23088.8575
$ws.Range("H110").Value = 5504.4
$ws.Range("I110").Value = 3508.3333
$ws.Range("J110").Value = 8498.5
$ws.Range("K110").Value = 3508.3333
$ws.Range("L110").Value = 8498.5
$ws.Range("M110").Value = -1463.3333
$ws.Range("N110").Value = -12588.5
$ws.Range("H136").Value = 10731.4375
$ws.Range("I136").Value = 11313.6
$ws.Range("K136").Value = 33940.8
$ws.Range("M136").Value = -31390.8
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1142530.9
$ws.Range("I94").Value = 1246215.5
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1246215.5
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -1245764.5
$ws.Range("N94").Value = -2902
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25006536
$ws.Range("I31").Value = 100001560
$ws.Range("K31").Value = 100001560
$ws.Range("M31").Value = -100001265
$ws.Range("H34").Value = 25006536
$ws.Range("I34").Value = 100001560
$ws.Range("K34").Value = 100001560
$ws.Range("M34").Value = -100001358
$ws.Range("H86").Value = 5267.5454
$ws.Range("I86").Value = 4749.75
$ws.Range("K86").Value = 4749.75
$ws.Range("M86").Value = -3626.75
$ws.Range("H89").Value = 5267.5454
$ws.Range("I89").Value = 4749.75
$ws.Range("K89").Value = 23748.75
$ws.Range("M89").Value = -18132.75
$ws.Range("H134").Value = 3682.5454
$ws.Range("I134").Value = 3780.8
$ws.Range("K134").Value = 11342.4
$ws.Range("M134").Value = -8807.400000000001
$ws.Range("H141").Value = 102050
$ws.Range("J141").Value = 102050
$ws.Range("L141").Value = 102050
$ws.Range("N141").Value = -112410
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 13550
$ws.Range("J63").Value = 25000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76498
$ws.Range("H66").Value = 13550
$ws.Range("J66").Value = 25000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -232488
$ws.Range("H80").Value = 5992.7144
$ws.Range("J80").Value = 5992.1665
$ws.Range("L80").Value = 17976.4995
$ws.Range("N80").Value = -19848.4995
$ws.Range("H83").Value = 5992.7144
$ws.Range("J83").Value = 5992.1665
$ws.Range("L83").Value = 53929.4985
$ws.Range("N83").Value = -63289.4985
$ws.Range("H107").Value = 2200.2144
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 2307.923
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 6923.768999999999
$ws.Range("M107").Value = -480
$ws.Range("N107").Value = -10763.769
$ws.Range("H114").Value = 2835
$ws.Range("J114").Value = 3321.5833
$ws.Range("L114").Value = 9964.749899999999
$ws.Range("N114").Value = -16472.7499
$ws.Range("H120").Value = 30785.715
$ws.Range("I120").Value = 5500
$ws.Range("J120").Value = 35000
$ws.Range("K120").Value = 16500
$ws.Range("L120").Value = 105000
$ws.Range("M120").Value = -11662
$ws.Range("N120").Value = -114676
$ws.Range("H129").Value = 3500
$ws.Range("I129").Value = 3500
$ws.Range("K129").Value = 10500
$ws.Range("M129").Value = -5500
$ws.Range("H131").Value = 8199267.5
$ws.Range("J131").Value = 6669298.5
$ws.Range("L131").Value = 20007895.5
$ws.Range("N131").Value = -20017975.5
$ws.Range("H139").Value = 478720.34
$ws.Range("I139").Value = 528190.9399999999
$ws.Range("J139").Value = 8749.5
$ws.Range("K139").Value = 1584572.82
$ws.Range("L139").Value = 26248.5
$ws.Range("M139").Value = -1579432.82
$ws.Range("N139").Value = -36528.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3052.1843
$ws.Range("I132").Value = 2385.1724
$ws.Range("J132").Value = 5201.4443
$ws.Range("K132").Value = 7155.5172
$ws.Range("L132").Value = 15604.3329
$ws.Range("M132").Value = -4625.5172
$ws.Range("N132").Value = -20664.3329
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1113.5238
$ws.Range("I22").Value = 999.7857
$ws.Range("K22").Value = 999.7857
$ws.Range("M22").Value = -704.7857
$ws.Range("H27").Value = 1113.5238
$ws.Range("I27").Value = 999.7857
$ws.Range("K27").Value = 999.7857
$ws.Range("M27").Value = -892.7857
$ws.Range("H46").Value = 6737.2383
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 6737.2383
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 6737.2383
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -7113.2383
$ws.Range("H61").Value = 4996.75
$ws.Range("I61").Value = 4998
$ws.Range("J61").Value = 4996.3335
$ws.Range("K61").Value = 4998
$ws.Range("L61").Value = 4996.3335
$ws.Range("M61").Value = -4796
$ws.Range("N61").Value = -5400.3335
$ws.Range("H113").Value = 4996.75
$ws.Range("I113").Value = 4998
$ws.Range("J113").Value = 4996.3335
$ws.Range("K113").Value = 4998
$ws.Range("L113").Value = 4996.3335
$ws.Range("M113").Value = -2828
$ws.Range("N113").Value = -9336.333500000001
$ws.Range("H122").Value = 43963828
$ws.Range("I122").Value = 90913060
$ws.Range("K122").Value = 272739180
$ws.Range("M122").Value = -272736730
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 959489.75
$ws.Range("I100").Value = 1568274.9
$ws.Range("K100").Value = 3136549.8
$ws.Range("M100").Value = -3136008.8
$ws.Range("H126").Value = 3458.4546
$ws.Range("I126").Value = 3387
$ws.Range("J126").Value = 3649
$ws.Range("K126").Value = 10161
$ws.Range("L126").Value = 10947
$ws.Range("M126").Value = -7691
$ws.Range("N126").Value = -15887
$ws.Range("H132").Value = 18522860
$ws.Range("I132").Value = 2416816.2
$ws.Range("J132").Value = 71442710
$ws.Range("K132").Value = 7250448.600000001
$ws.Range("L132").Value = 214328130
$ws.Range("M132").Value = -7247918.600000001
$ws.Range("N132").Value = -214333190
$ws.Range("H136").Value = 9076.352000000001
$ws.Range("I136").Value = 5088.5
$ws.Range("J136").Value = 11795.341
$ws.Range("K136").Value = 15265.5
$ws.Range("L136").Value = 35386.023
$ws.Range("M136").Value = -12715.5
$ws.Range("N136").Value = -40486.023
